# Elle giriş ve hata çözümlemesi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 3 (Ikigai / İnsan / Atomik occupy
# rows 2-4 today) so there is room for the new "Scikit-Learn" entry while the
# other rows' relative order stays intact.
$ws.Rows.Item(3).Insert()

# Row 2: "Atomik Aliskanliklar" (previously row 4) moves up to the top data row.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Atomik Aliskanliklar - Kücük Degisikler Büyük Sonuclar - Kötü Aliskanliklardan Kurtulup Iyi Aliskanliklar Edinmek Icin Kolay ve Etkisi Kanitlanmis Bir Yöntem"
$ws.Range("C2").Value = "['James Clear']"

# Row 3 (new row): "Scikit-Learn Makine Ogrenimi" by Oreilly.
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Scikit-Learn Makine Ogrenimi"
$ws.Range("C3").Value = "Oreilly"

# Row 4: "Ikigai" (previously row 2), renumbered.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Ikigai - Uygulama Rehberi - Japonlarin Uzun ve Mutlu Yasam Sirrini Hayata Gecirin"
$ws.Range("C4").Value = "['Hector Garcia', 'Francesc Miralles']"

# Row 5: "İnsan tabiatını tanıma" (previously row 3), renumbered.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "İnsan tabiatını tanıma"
$ws.Range("C5").Value = "['Alfred Adler']"
